$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $startsWith) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.StartsWith($startsWith)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. "...installed and winscp please follow..." -> "...installed and WinSCP
#    please follow..." (capitalize the product name).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "If you do not have the tomcat component installed and winscp please follow these lines, if you do have them please skip the first three points.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If you do not have the tomcat component installed and WinSCP please follow these lines, if you do have them please skip the first three points.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Replace the "Copy tomcat8.tar, winscp577setup.exe and putty.exe ..."
#    paragraph with "Install WinSCP on your local machine: " followed by a
#    real hyperlink to the WinSCP download page.
# ---------------------------------------------------------------------------
$copyIdx = Find-ParagraphIndex $d "Copy tomcat8.tar"
$copyPara = $d.Paragraphs($copyIdx)
$copyRange = $copyPara.Range
$copyRange.End = $copyRange.End - 1
$copyRange.Text = "Install WinSCP on your local machine: https://winscp.net/eng/download.php "

$paraText = $copyPara.Range.Text
$urlStart = $copyPara.Range.Start + $paraText.IndexOf("https://winscp.net/eng/download.php")
$urlEnd = $urlStart + "https://winscp.net/eng/download.php".Length
$urlRange = $d.Range($urlStart, $urlEnd)
$d.Hyperlinks.Add($urlRange, "https://winscp.net/eng/download.php") | Out-Null

# ---------------------------------------------------------------------------
# 3. Remove the paragraphs describing the old docker/tomcat copy workflow:
#    "Open a command prompt ...", "docker load -i tomcat8.tar",
#    "docker run --name tomcat ...", "docker start zth/tomcat",
#    an empty paragraph, "Install winscp on you local machine." and another
#    empty paragraph. These steps are no longer needed.
# ---------------------------------------------------------------------------
$startIdx = Find-ParagraphIndex $d "Open a command prompt"
$endIdx = Find-ParagraphIndex $d "Install winscp on you local machine."
$startPara = $d.Paragraphs($startIdx)
$endPara = $d.Paragraphs($endIdx + 1)
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. Reword "Open winscp and connect to the tomcat server on docker like in
#    the printscreen ..." -> "Open WinSCP and connect to the tomcat server on
#    Docker like in the screenshot ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Open winscp and connect to the tomcat server on docker like in the printscreen (the username and password are root and root):",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Open WinSCP and connect to the tomcat server on Docker like in the screenshot (the username and password are root and root):",
    2) | Out-Null

Write-Host "Edit complete"
